# PhanCongCongViec.xlsx — "Hoàn thành tất cả. Thêm bảng phân công"
#
# The task table (rows 10-14, under "Đặng Công Thắng") was empty before;
# this fills in the 5 task descriptions and marks every one of them as
# 100% complete (a plain "1" formatted as a percentage, matching the
# existing C3 cell's look).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: column B needs to fit the long task descriptions,
# and the new column C (progress %) gets its own width too. ---
$ws.Columns.Item(2).ColumnWidth = 41.43
$ws.Columns.Item(3).ColumnWidth = 13.8

# --- Task descriptions (column B), rows 10-14 ---
$ws.Range("B10").Value = "1. Hoàn thành layout login, register"
$ws.Range("B11").Value = "2. Liên kết API axios, localStrorage"
$ws.Range("B13").Value = "4. Hiển thị Name ra trang index"
$ws.Range("B12").Value = "3. Hoàn thành sự kiện đăng nhập,  đăng xuất"
$ws.Range("B14").Value = "5. Hoàn thành kiểm tra lỗi"

# --- Progress (column C): all five tasks done = 100% ---
$ws.Range("C10:C14").Value = 1
$ws.Range("C10:C14").NumberFormat = "0%"

# --- Leave the cursor where the author last left it ---
$ws.Range("H12").Select()
